# Insert a new data row at row 19 (pushes existing rows 19-53 down to 20-54)
# and populate it with a new weekly price observation, consistent with the
# rest of the "Arveja Verde" series for Vega Monumental Concepción.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Insert()

$newRow = 19

$ws.Cells.Item($newRow, 1).Value = 11
$ws.Cells.Item($newRow, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($newRow, 3).Value = "Bíobío"
$ws.Cells.Item($newRow, 4).Value2 = 45259
$ws.Cells.Item($newRow, 5).Value = 8
$ws.Cells.Item($newRow, 6).Value = 100112022
$ws.Cells.Item($newRow, 7).Value = "Arveja Verde"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 180
$ws.Cells.Item($newRow, 11).Value = 18000
$ws.Cells.Item($newRow, 12).Value = 20000
$ws.Cells.Item($newRow, 13).Value = 18889
$ws.Cells.Item($newRow, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item($newRow, 15).Value = "Región del Maule"
$ws.Cells.Item($newRow, 16).Value = 756
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
